# Append 5 new daily-report rows (234-238) to the bottom of the table,
# continuing the date series 44308-44312 (22-26 Apr 2021), matching the
# formatting of the preceding data row (A column keeps the date style/
# number format, B/C/D stay unstyled numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 234; A = 44308; B = 2; C = 16; D = 106.0515675747332 },
    @{ Row = 235; A = 44309; B = 0; C = 14; D = 92.79512162789156 },
    @{ Row = 236; A = 44310; B = 0; C = 7;  D = 46.39756081394578 },
    @{ Row = 237; A = 44311; B = 3; C = 6;  D = 39.76933784052495 },
    @{ Row = 238; A = 44312; B = 2; C = 7;  D = 46.39756081394578 }
)

# Use the last existing row (233) as the formatting template for the new
# A-column cells (it already carries the date number format + style).
$templateCell = $ws.Range("A233")
$templateCell.Copy() | Out-Null

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D

    # xlPasteFormats = -4122 : copy only the formatting from A233 onto the
    # new date cell so it keeps the same style (centered, bordered, date fmt).
    $ws.Range("A$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
